$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = "Codificadió básica de l'objecte Jugador"
$ws.Range("A31").Select() | Out-Null
$ws.Range("B31").Value = "Planificada"
$ws.Range("B31").Select() | Out-Null
$ws.Range("C31").Value = 1
$ws.Range("C31").Select() | Out-Null
$ws.Range("D31").Value = "Aleix"
$ws.Range("D31").Select() | Out-Null
$ws.Range("E31").Value = 16
$ws.Range("E31").Select() | Out-Null
$ws.Range("F31").Value = "Entitat que guarda la puntuació, diners, llista de totems i proveeix metodes per actualitzar i recuperar aquests atributs."
$ws.Range("F31").Select() | Out-Null
